$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1874.9131
$ws.Range("J19").Value = 2017.25
$ws.Range("L19").Value = 2017.25
$ws.Range("N19").Value = -2367.25

$ws.Range("H86").Value = 3279.9
$ws.Range("I86").Value = 3272.375
$ws.Range("J86").Value = 3284.9167
$ws.Range("K86").Value = 3272.375
$ws.Range("L86").Value = 3284.9167
$ws.Range("M86").Value = -2149.375
$ws.Range("N86").Value = -5530.9167

$ws.Range("H89").Value = 3279.9
$ws.Range("I89").Value = 3272.375
$ws.Range("J89").Value = 3284.9167
$ws.Range("K89").Value = 16361.875
$ws.Range("L89").Value = 16424.5835
$ws.Range("M89").Value = -10745.875
$ws.Range("N89").Value = -27656.5835

$ws.Range("H100").Value = 1854.4
$ws.Range("I100").Value = 1681.1428
$ws.Range("K100").Value = 1681.1428
$ws.Range("M100").Value = -1140.1428

$ws.Range("H103").Value = 821.8889
$ws.Range("I103").Value = 583
$ws.Range("J103").Value = 1299.6666
$ws.Range("K103").Value = 1749
$ws.Range("L103").Value = 3898.9998
$ws.Range("M103").Value = -1163
$ws.Range("N103").Value = -5070.9998

$ws.Range("H116").Value = 8979.286
$ws.Range("J116").Value = 7750.1816
$ws.Range("L116").Value = 7750.1816
$ws.Range("N116").Value = -14634.1816

$ws.Range("H137").Value = 5139.5557
$ws.Range("I137").Value = 4510.2
$ws.Range("K137").Value = 13530.6
$ws.Range("M137").Value = -10980.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2023.23
$ws.Range("I32").Value = 1288.409
$ws.Range("K32").Value = 1288.409
$ws.Range("M32").Value = -1001.409

$ws.Range("H45").Value = 7940.6
$ws.Range("J45").Value = 2329
$ws.Range("L45").Value = 2329
$ws.Range("N45").Value = -3083

$ws.Range("H74").Value = 1819.6875
$ws.Range("I74").Value = 1466.75
$ws.Range("K74").Value = 1466.75
$ws.Range("M74").Value = -592.75

$ws.Range("H77").Value = 1819.6875
$ws.Range("I77").Value = 1466.75
$ws.Range("K77").Value = 7333.75
$ws.Range("M77").Value = -2965.75

$ws.Range("H102").Value = 4784.05
$ws.Range("J102").Value = 4884.625
$ws.Range("L102").Value = 4884.625
$ws.Range("N102").Value = -8128.625

$ws.Range("H110").Value = 1218.6
$ws.Range("I110").Value = 1218.6
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1218.6
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 826.4000000000001
$ws.Range("N110").ClearContents()

$ws.Range("H114").Value = 70400
$ws.Range("J114").Value = 70400
$ws.Range("L114").Value = 70400
$ws.Range("N114").Value = -79078

$ws.Range("H122").Value = 5919.729
$ws.Range("I122").Value = 3335.7715
$ws.Range("J122").Value = 7402.3276
$ws.Range("K122").Value = 10007.3145
$ws.Range("L122").Value = 22206.9828
$ws.Range("M122").Value = -7557.3145
$ws.Range("N122").Value = -27106.9828

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 286.66666
$ws.Range("I22").Value = 154.5
$ws.Range("J22").Value = 551
$ws.Range("K22").Value = 154.5
$ws.Range("L22").Value = 551
$ws.Range("M22").Value = 18.5
$ws.Range("N22").Value = -897

$ws.Range("H86").Value = 3390.158
$ws.Range("I86").Value = 3061.1428
$ws.Range("J86").Value = 3582.0833
$ws.Range("K86").Value = 3061.1428
$ws.Range("L86").Value = 3582.0833
$ws.Range("M86").Value = -1938.1428
$ws.Range("N86").Value = -5828.0833

$ws.Range("H89").Value = 3390.158
$ws.Range("I89").Value = 3061.1428
$ws.Range("J89").Value = 3582.0833
$ws.Range("K89").Value = 15305.714
$ws.Range("L89").Value = 17910.4165
$ws.Range("M89").Value = -9689.714
$ws.Range("N89").Value = -29142.4165

$ws.Range("H99").Value = 5315.1665
$ws.Range("I99").Value = 5998.643
$ws.Range("K99").Value = 5998.643
$ws.Range("M99").Value = -4500.643

$ws.Range("H105").Value = 2064.6155
$ws.Range("I105").Value = 2069.4255
$ws.Range("K105").Value = 2069.4255
$ws.Range("M105").Value = -322.4254999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1350.8928
$ws.Range("I16").Value = 951.375
$ws.Range("J16").Value = 1510.7
$ws.Range("K16").Value = 951.375
$ws.Range("L16").Value = 1510.7
$ws.Range("M16").Value = -664.375
$ws.Range("N16").Value = -2084.7

$ws.Range("H58").Value = 11096.407
$ws.Range("I58").Value = 10782.333
$ws.Range("J58").Value = 12195.667
$ws.Range("K58").Value = 10782.333
$ws.Range("L58").Value = 12195.667
$ws.Range("M58").Value = -10579.333
$ws.Range("N58").Value = -12601.667

$ws.Range("H105").Value = 1247.4117
$ws.Range("I105").Value = 1247.4117
$ws.Range("K105").Value = 1247.4117
$ws.Range("M105").Value = 499.5882999999999

$ws.Range("H113").Value = 1350.8928
$ws.Range("I113").Value = 951.375
$ws.Range("J113").Value = 1510.7
$ws.Range("K113").Value = 951.375
$ws.Range("L113").Value = 1510.7
$ws.Range("M113").Value = 1218.625
$ws.Range("N113").Value = -5850.7

$ws.Range("H136").Value = 11096.407
$ws.Range("I136").Value = 10782.333
$ws.Range("J136").Value = 12195.667
$ws.Range("K136").Value = 32346.999
$ws.Range("L136").Value = 36587.001
$ws.Range("M136").Value = -29796.999
$ws.Range("N136").Value = -41687.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3891.125
$ws.Range("I80").Value = 3777.3333
$ws.Range("J80").Value = 3959.4
$ws.Range("K80").Value = 11331.9999
$ws.Range("L80").Value = 11878.2
$ws.Range("M80").Value = -10395.9999
$ws.Range("N80").Value = -13750.2

$ws.Range("H83").Value = 3891.125
$ws.Range("I83").Value = 3777.3333
$ws.Range("J83").Value = 3959.4
$ws.Range("K83").Value = 33995.9997
$ws.Range("L83").Value = 35634.6
$ws.Range("M83").Value = -29315.9997
$ws.Range("N83").Value = -44994.6

$ws.Range("H129").Value = 13368389
$ws.Range("J129").Value = 13891643
$ws.Range("L129").Value = 41674929
$ws.Range("N129").Value = -41684929

$ws.Range("H131").Value = 4368521
$ws.Range("I131").Value = 24064746
$ws.Range("J131").Value = 2316830.8
$ws.Range("K131").Value = 72194238
$ws.Range("L131").Value = 6950492.399999999
$ws.Range("M131").Value = -72189198
$ws.Range("N131").Value = -6960572.399999999

$ws.Range("H132").Value = 1630.2858
$ws.Range("I132").Value = 1081.3334
$ws.Range("K132").Value = 9732.000599999999
$ws.Range("M132").Value = -7202.000599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H102").Value = 4939.6113
$ws.Range("I102").Value = 5761.2593
$ws.Range("K102").Value = 5761.2593
$ws.Range("M102").Value = -4139.2593

$ws.Range("H113").Value = 4515.846
$ws.Range("I113").Value = 4275.2856
$ws.Range("J113").Value = 4796.5
$ws.Range("K113").Value = 4275.2856
$ws.Range("L113").Value = 4796.5
$ws.Range("M113").Value = -2105.2856
$ws.Range("N113").Value = -9136.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680

$ws.Range("H122").Value = 3490.1177
$ws.Range("I122").Value = 3331
$ws.Range("K122").Value = 9993
$ws.Range("M122").Value = -7543

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 14999.5
$ws.Range("I22").Value = 13999
$ws.Range("J22").Value = 16000
$ws.Range("K22").Value = 13999
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = -13706
$ws.Range("N22").Value = -16586

$ws.Range("H24").Value = 19333.777
$ws.Range("J24").Value = 19333.777
$ws.Range("L24").Value = 19333.777
$ws.Range("N24").Value = -19793.777

$ws.Range("H81").Value = 1217102.9
$ws.Range("I81").Value = 12880.667
$ws.Range("J81").Value = 3384702.8
$ws.Range("K81").Value = 25761.334
$ws.Range("L81").Value = 6769405.6
$ws.Range("M81").Value = -24700.334
$ws.Range("N81").Value = -6771527.6

$ws.Range("H84").Value = 1217102.9
$ws.Range("I84").Value = 12880.667
$ws.Range("J84").Value = 3384702.8
$ws.Range("K84").Value = 128806.67
$ws.Range("L84").Value = 33847028
$ws.Range("M84").Value = -123502.67
$ws.Range("N84").Value = -33857636

$ws.Range("H110").Value = 34500
$ws.Range("J110").Value = 34500
$ws.Range("L110").Value = 34500
$ws.Range("N110").Value = -42680
